# Apply odds/score updates for Jogos_da_Semana_FlashScore_2024-10-29.xlsx
# (values taken from FlashScore refresh of rows 2-5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 7.5
$ws.Range("Q2").Value = 1.48
$ws.Range("R2").Value = 2.6
$ws.Range("U2").Value = 1.67
$ws.Range("V2").Value = 2.1
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 10
$ws.Range("AU2").Value = 8
$ws.Range("AV2").Value = 41
$ws.Range("AW2").Value = 8.5
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 101
$ws.Range("BB2").Value = 151
# Row 3
$ws.Range("G3").Value = 1.95
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 4.2
$ws.Range("J3").Value = 2.75
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 8
$ws.Range("Z3").Value = 17
$ws.Range("AA3").Value = 19
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 7
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 67
$ws.Range("AH3").Value = 9
$ws.Range("AI3").Value = 19
$ws.Range("AJ3").Value = 15
$ws.Range("AK3").Value = 41
$ws.Range("AL3").Value = 41
$ws.Range("AN3").Value = 3.75
$ws.Range("AO3").Value = 11
$ws.Range("AQ3").Value = 41
$ws.Range("AU3").Value = 9
$ws.Range("AV3").Value = 67
$ws.Range("AW3").Value = 6
$ws.Range("AX3").Value = 23
$ws.Range("AY3").Value = 34
$ws.Range("AZ3").Value = 81
$ws.Range("BA3").Value = 126
$ws.Range("BB3").Value = 351
# Row 4
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
# Row 5
$ws.Range("K5").Value = 2.1
$ws.Range("L5").Value = 6.5
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 2.2
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("W5").Value = 5.5
$ws.Range("X5").Value = 6.5
$ws.Range("AB5").Value = 34
$ws.Range("AC5").Value = 8
$ws.Range("AE5").Value = 21
$ws.Range("AF5").Value = 81
$ws.Range("AH5").Value = 13
$ws.Range("AM5").Value = 51
$ws.Range("AQ5").Value = 26
$ws.Range("AS5").Value = 201
$ws.Range("AT5").Value = 2.63
$ws.Range("BA5").Value = 201
